# Fix Excel copy/paste error:
# The "Term Accession Number (DPBO:1000161)" column (D) on the
# "Chromatography" sheet should contain the same term accession number
# (http://purl.obolibrary.org/obo/NCIT_C16431) for every row, matching
# column B's constant "Chromatography" Protocol Type. Rows 3-5 had been
# accidentally overwritten with unrelated NCIT accession numbers during a
# copy/paste - restore the correct, consistent value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chromatography")

$ws.Range("D3").Value = "http://purl.obolibrary.org/obo/NCIT_C16431"
$ws.Range("D4").Value = "http://purl.obolibrary.org/obo/NCIT_C16431"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/NCIT_C16431"

# While fixing the error, the hidden helper columns C and D were revealed
# so the mistake could be spotted and corrected.
$ws.Columns.Item(3).Hidden = $false
$ws.Columns.Item(4).Hidden = $false

# The Chromatography sheet became the active sheet/tab, with D6 selected
# as the active cell after the correction.
$ws.Activate()
$ws.Range("D6").Select()
